# Auto-generated edit script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# across several leve-profit worksheets, per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 38508680
$ws.Cells.Item(33, 9).Value = 50001016
$ws.Cells.Item(33, 10).Value = 200900.67
$ws.Cells.Item(33, 11).Value = 50001016
$ws.Cells.Item(33, 12).Value = 200900.67
$ws.Cells.Item(33, 13).Value = -50000787
$ws.Cells.Item(33, 14).Value = -201358.67
# Row 43
$ws.Cells.Item(43, 8).Value = 1115.4166
$ws.Cells.Item(43, 10).Value = 750.6667
$ws.Cells.Item(43, 12).Value = 750.6667
$ws.Cells.Item(43, 14).Value = -888.6667
# Row 53
$ws.Cells.Item(53, 8).Value = 149.36363
$ws.Cells.Item(53, 9).Value = 177.57143
$ws.Cells.Item(53, 11).Value = 177.57143
$ws.Cells.Item(53, 13).Value = 459.42857
# Row 64
$ws.Cells.Item(64, 8).Value = 3175.7188
$ws.Cells.Item(64, 9).Value = 2831.6667
$ws.Cells.Item(64, 10).Value = 3618.0715
$ws.Cells.Item(64, 11).Value = 2831.6667
$ws.Cells.Item(64, 12).Value = 3618.0715
$ws.Cells.Item(64, 13).Value = -2583.6667
$ws.Cells.Item(64, 14).Value = -4114.0715
# Row 67
$ws.Cells.Item(67, 8).Value = 3175.7188
$ws.Cells.Item(67, 9).Value = 2831.6667
$ws.Cells.Item(67, 10).Value = 3618.0715
$ws.Cells.Item(67, 11).Value = 2831.6667
$ws.Cells.Item(67, 12).Value = 3618.0715
$ws.Cells.Item(67, 13).Value = -1973.6667
$ws.Cells.Item(67, 14).Value = -5334.0715
# Row 76
$ws.Cells.Item(76, 8).Value = 22434.8
$ws.Cells.Item(76, 9).Value = 2990
$ws.Cells.Item(76, 10).Value = 27296
$ws.Cells.Item(76, 11).Value = 2990
$ws.Cells.Item(76, 12).Value = 27296
$ws.Cells.Item(76, 13).Value = -2675
$ws.Cells.Item(76, 14).Value = -27926
# Row 79
$ws.Cells.Item(79, 8).Value = 22434.8
$ws.Cells.Item(79, 9).Value = 2990
$ws.Cells.Item(79, 10).Value = 27296
$ws.Cells.Item(79, 11).Value = 2990
$ws.Cells.Item(79, 12).Value = 27296
$ws.Cells.Item(79, 13).Value = -1898
$ws.Cells.Item(79, 14).Value = -29480
# Row 86
$ws.Cells.Item(86, 8).Value = 32629.572
$ws.Cells.Item(86, 9).Value = 60000
$ws.Cells.Item(86, 10).Value = 21681.4
$ws.Cells.Item(86, 11).Value = 60000
$ws.Cells.Item(86, 12).Value = 21681.4
$ws.Cells.Item(86, 13).Value = -58877
$ws.Cells.Item(86, 14).Value = -23927.4
# Row 89
$ws.Cells.Item(89, 8).Value = 32629.572
$ws.Cells.Item(89, 9).Value = 60000
$ws.Cells.Item(89, 10).Value = 21681.4
$ws.Cells.Item(89, 11).Value = 300000
$ws.Cells.Item(89, 12).Value = 108407
$ws.Cells.Item(89, 13).Value = -294384
$ws.Cells.Item(89, 14).Value = -119639
# Row 98
$ws.Cells.Item(98, 8).Value = 1658.5
$ws.Cells.Item(98, 9).Value = 1478.3846
$ws.Cells.Item(98, 10).Value = 4000
$ws.Cells.Item(98, 11).Value = 1478.3846
$ws.Cells.Item(98, 12).Value = 4000
$ws.Cells.Item(98, 13).Value = 19.61539999999991
$ws.Cells.Item(98, 14).Value = -6996
# Row 106
$ws.Cells.Item(106, 8).Value = 1700
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 1700
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 1700
$ws.Cells.Item(106, 13).ClearContents() | Out-Null
$ws.Cells.Item(106, 14).Value = -2962
# Row 116
$ws.Cells.Item(116, 8).Value = 4543.92
$ws.Cells.Item(116, 9).Value = 2285.4285
$ws.Cells.Item(116, 11).Value = 2285.4285
$ws.Cells.Item(116, 13).Value = 1156.5715
# Row 122
$ws.Cells.Item(122, 8).Value = 1658.5
$ws.Cells.Item(122, 9).Value = 1478.3846
$ws.Cells.Item(122, 10).Value = 4000
$ws.Cells.Item(122, 11).Value = 4435.1538
$ws.Cells.Item(122, 12).Value = 12000
$ws.Cells.Item(122, 13).Value = -1985.1538
$ws.Cells.Item(122, 14).Value = -16900
# Row 123
$ws.Cells.Item(123, 8).Value = 29780
$ws.Cells.Item(123, 10).Value = 29780
$ws.Cells.Item(123, 12).Value = 29780
$ws.Cells.Item(123, 14).Value = -39580
# Row 127
$ws.Cells.Item(127, 8).Value = 1236.4
$ws.Cells.Item(127, 9).Value = 703
$ws.Cells.Item(127, 10).Value = 1655.5
$ws.Cells.Item(127, 11).Value = 2109
$ws.Cells.Item(127, 12).Value = 4966.5
$ws.Cells.Item(127, 13).Value = 2851
$ws.Cells.Item(127, 14).Value = -14886.5
# Row 132
$ws.Cells.Item(132, 8).Value = 1943.5883
$ws.Cells.Item(132, 9).Value = 1227.6154
$ws.Cells.Item(132, 10).Value = 4270.5
$ws.Cells.Item(132, 11).Value = 3682.8462
$ws.Cells.Item(132, 12).Value = 12811.5
$ws.Cells.Item(132, 13).Value = -1152.8462
$ws.Cells.Item(132, 14).Value = -17871.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 12531
$ws.Cells.Item(32, 9).Value = 4277.5386
$ws.Cells.Item(32, 11).Value = 4277.5386
$ws.Cells.Item(32, 13).Value = -3990.5386
# Row 61
$ws.Cells.Item(61, 8).Value = 1960.1818
$ws.Cells.Item(61, 9).Value = 1755.5454
$ws.Cells.Item(61, 10).Value = 2164.818
$ws.Cells.Item(61, 11).Value = 1755.5454
$ws.Cells.Item(61, 12).Value = 2164.818
$ws.Cells.Item(61, 13).Value = -1543.5454
$ws.Cells.Item(61, 14).Value = -2588.818
# Row 74
$ws.Cells.Item(74, 8).Value = 919.5
$ws.Cells.Item(74, 9).Value = 690.37933
$ws.Cells.Item(74, 10).Value = 1310.3529
$ws.Cells.Item(74, 11).Value = 690.37933
$ws.Cells.Item(74, 12).Value = 1310.3529
$ws.Cells.Item(74, 13).Value = 183.62067
$ws.Cells.Item(74, 14).Value = -3058.3529
# Row 77
$ws.Cells.Item(77, 8).Value = 919.5
$ws.Cells.Item(77, 9).Value = 690.37933
$ws.Cells.Item(77, 10).Value = 1310.3529
$ws.Cells.Item(77, 11).Value = 3451.89665
$ws.Cells.Item(77, 12).Value = 6551.7645
$ws.Cells.Item(77, 13).Value = 916.1033500000003
$ws.Cells.Item(77, 14).Value = -15287.7645
# Row 122
$ws.Cells.Item(122, 8).Value = 1850.7097
$ws.Cells.Item(122, 9).Value = 1876.8422
$ws.Cells.Item(122, 10).Value = 1809.3334
$ws.Cells.Item(122, 11).Value = 5630.5266
$ws.Cells.Item(122, 12).Value = 5428.0002
$ws.Cells.Item(122, 13).Value = -3180.5266
$ws.Cells.Item(122, 14).Value = -10328.0002
# Row 136
$ws.Cells.Item(136, 8).Value = 1960.1818
$ws.Cells.Item(136, 9).Value = 1755.5454
$ws.Cells.Item(136, 10).Value = 2164.818
$ws.Cells.Item(136, 11).Value = 5266.6362
$ws.Cells.Item(136, 12).Value = 6494.454000000001
$ws.Cells.Item(136, 13).Value = -2716.6362
$ws.Cells.Item(136, 14).Value = -11594.454

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2260.3845
$ws.Cells.Item(31, 9).Value = 1530.5264
$ws.Cells.Item(31, 10).Value = 4241.4287
$ws.Cells.Item(31, 11).Value = 1530.5264
$ws.Cells.Item(31, 12).Value = 4241.4287
$ws.Cells.Item(31, 13).Value = -1235.5264
$ws.Cells.Item(31, 14).Value = -4831.4287
# Row 34
$ws.Cells.Item(34, 8).Value = 2260.3845
$ws.Cells.Item(34, 9).Value = 1530.5264
$ws.Cells.Item(34, 10).Value = 4241.4287
$ws.Cells.Item(34, 11).Value = 1530.5264
$ws.Cells.Item(34, 12).Value = 4241.4287
$ws.Cells.Item(34, 13).Value = -1328.5264
$ws.Cells.Item(34, 14).Value = -4645.4287
# Row 86
$ws.Cells.Item(86, 8).Value = 24339.652
$ws.Cells.Item(86, 9).Value = 37571.715
$ws.Cells.Item(86, 10).Value = 3756.4443
$ws.Cells.Item(86, 11).Value = 37571.715
$ws.Cells.Item(86, 12).Value = 3756.4443
$ws.Cells.Item(86, 13).Value = -36448.715
$ws.Cells.Item(86, 14).Value = -6002.4443
# Row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents() | Out-Null
# Row 89
$ws.Cells.Item(89, 8).Value = 24339.652
$ws.Cells.Item(89, 9).Value = 37571.715
$ws.Cells.Item(89, 10).Value = 3756.4443
$ws.Cells.Item(89, 11).Value = 187858.575
$ws.Cells.Item(89, 12).Value = 18782.2215
$ws.Cells.Item(89, 13).Value = -182242.575
$ws.Cells.Item(89, 14).Value = -30014.2215
# Row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents() | Out-Null
# Row 99
$ws.Cells.Item(99, 8).Value = 1770.25
$ws.Cells.Item(99, 9).Value = 1408
$ws.Cells.Item(99, 10).Value = 2857
$ws.Cells.Item(99, 11).Value = 1408
$ws.Cells.Item(99, 12).Value = 2857
$ws.Cells.Item(99, 13).Value = 90
$ws.Cells.Item(99, 14).Value = -5853
# Row 126
$ws.Cells.Item(126, 8).Value = 1770.25
$ws.Cells.Item(126, 9).Value = 1408
$ws.Cells.Item(126, 10).Value = 2857
$ws.Cells.Item(126, 11).Value = 4224
$ws.Cells.Item(126, 12).Value = 8571
$ws.Cells.Item(126, 13).Value = -1754
$ws.Cells.Item(126, 14).Value = -13511
# Row 134
$ws.Cells.Item(134, 8).Value = 1671.1887
$ws.Cells.Item(134, 9).Value = 1362.1143
$ws.Cells.Item(134, 10).Value = 2272.1667
$ws.Cells.Item(134, 11).Value = 4086.3429
$ws.Cells.Item(134, 12).Value = 6816.500100000001
$ws.Cells.Item(134, 13).Value = -1551.3429
$ws.Cells.Item(134, 14).Value = -11886.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 2661.7354
$ws.Cells.Item(102, 9).Value = 2259.2964
$ws.Cells.Item(102, 10).Value = 4214
$ws.Cells.Item(102, 11).Value = 2259.2964
$ws.Cells.Item(102, 12).Value = 4214
$ws.Cells.Item(102, 13).Value = -637.2964000000002
$ws.Cells.Item(102, 14).Value = -7458
# Row 122
$ws.Cells.Item(122, 8).Value = 1722.5416
$ws.Cells.Item(122, 9).Value = 1540.7778
$ws.Cells.Item(122, 10).Value = 1831.6
$ws.Cells.Item(122, 11).Value = 4622.3334
$ws.Cells.Item(122, 12).Value = 5494.799999999999
$ws.Cells.Item(122, 13).Value = -2172.3334
$ws.Cells.Item(122, 14).Value = -10394.8
# Row 126
$ws.Cells.Item(126, 8).Value = 1772.1538
$ws.Cells.Item(126, 9).Value = 1413.875
$ws.Cells.Item(126, 10).Value = 2345.4
$ws.Cells.Item(126, 11).Value = 4241.625
$ws.Cells.Item(126, 12).Value = 7036.200000000001
$ws.Cells.Item(126, 13).Value = -1771.625
$ws.Cells.Item(126, 14).Value = -11976.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 2199.75
$ws.Cells.Item(7, 9).Value = 1884.7
$ws.Cells.Item(7, 10).Value = 3775
$ws.Cells.Item(7, 11).Value = 1884.7
$ws.Cells.Item(7, 12).Value = 3775
$ws.Cells.Item(7, 13).Value = -1772.7
$ws.Cells.Item(7, 14).Value = -3999
# Row 40
$ws.Cells.Item(40, 8).Value = 2839.7837
$ws.Cells.Item(40, 9).Value = 2617.1785
$ws.Cells.Item(40, 10).Value = 3532.3333
$ws.Cells.Item(40, 11).Value = 2617.1785
$ws.Cells.Item(40, 12).Value = 3532.3333
$ws.Cells.Item(40, 13).Value = -2481.1785
$ws.Cells.Item(40, 14).Value = -3804.3333
# Row 46
$ws.Cells.Item(46, 8).Value = 1471.5714
$ws.Cells.Item(46, 9).Value = 1450.25
$ws.Cells.Item(46, 11).Value = 1450.25
$ws.Cells.Item(46, 13).Value = -1262.25
# Row 55
$ws.Cells.Item(55, 8).Value = 268.8125
$ws.Cells.Item(55, 10).Value = 345.7143
$ws.Cells.Item(55, 12).Value = 345.7143
$ws.Cells.Item(55, 14).Value = -691.7143
# Row 122
$ws.Cells.Item(122, 8).Value = 12793757
$ws.Cells.Item(122, 9).Value = 1669675.1
$ws.Cells.Item(122, 10).Value = 38464716
$ws.Cells.Item(122, 11).Value = 5009025.300000001
$ws.Cells.Item(122, 12).Value = 115394148
$ws.Cells.Item(122, 13).Value = -5006575.300000001
$ws.Cells.Item(122, 14).Value = -115399048
# Row 126
$ws.Cells.Item(126, 8).Value = 2199.75
$ws.Cells.Item(126, 9).Value = 1884.7
$ws.Cells.Item(126, 10).Value = 3775
$ws.Cells.Item(126, 11).Value = 5654.1
$ws.Cells.Item(126, 12).Value = 11325
$ws.Cells.Item(126, 13).Value = -3184.1
$ws.Cells.Item(126, 14).Value = -16265
# Row 132
$ws.Cells.Item(132, 8).Value = 3177.3076
$ws.Cells.Item(132, 9).Value = 2363.25
$ws.Cells.Item(132, 10).Value = 4479.8
$ws.Cells.Item(132, 11).Value = 7089.75
$ws.Cells.Item(132, 12).Value = 13439.4
$ws.Cells.Item(132, 13).Value = -4559.75
$ws.Cells.Item(132, 14).Value = -18499.4

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 3584.7942
$ws.Cells.Item(107, 9).Value = 720.381
$ws.Cells.Item(107, 10).Value = 8211.923000000001
$ws.Cells.Item(107, 11).Value = 2161.143
$ws.Cells.Item(107, 12).Value = 24635.769
$ws.Cells.Item(107, 13).Value = -241.143
$ws.Cells.Item(107, 14).Value = -28475.769
# Row 122
$ws.Cells.Item(122, 8).Value = 2303.5334
$ws.Cells.Item(122, 9).Value = 1752.5454
$ws.Cells.Item(122, 10).Value = 3818.75
$ws.Cells.Item(122, 11).Value = 5257.6362
$ws.Cells.Item(122, 12).Value = 11456.25
$ws.Cells.Item(122, 13).Value = -2807.6362
$ws.Cells.Item(122, 14).Value = -16356.25

Write-Host "Applied market-data refresh updates."